# Append the 04-12-2025 gold-price row (row 79) to Sheet1, matching the
# style/formatting already used by the preceding rows (A: bordered text,
# B: bordered + word-wrapped text) without disturbing the workbook's
# existing style table or introducing any date auto-conversion.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 79

$dateValue = "04-12-2025"
$priceValue = "The price of gold in India today is ₹13,036 per gram for 24 karat gold, ₹11,950 per gram for 22 karat gold and ₹9,778 per gram for 18 karat gold (also called 999 gold)."

# Enter the values as text-returning formulas first (so the date-like
# string "04-12-2025" isn't auto-parsed into a date serial number), then
# collapse them to plain text values in place. This preserves the same
# cell style used by row 78 instead of minting a brand-new number-format
# style for the new row.
$ws.Cells.Item($newRow, 1).Formula = '="' + $dateValue + '"'
$ws.Cells.Item($newRow, 2).Formula = '="' + $priceValue + '"'

$rowRange = "A" + $newRow + ":B" + $newRow
$ws.Range($rowRange).Copy() | Out-Null
$ws.Range($rowRange).PasteSpecial(-4163) | Out-Null
